# Applies the "Add Block Diagram Files" edit to the Functional Requirements
# document: a handful of grammar-checker paragraph-splits (w:proofErr
# markers), a consolidation of the "assistant" bullets into a single
# "forward complaints" bullet, a wording tweak, and moving the
# w:lastRenderedPageBreak marker from one run to another.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match $needle) {
            return $p
        }
    }
    return $null
}

function Set-ParaXml($para, $innerXml) {
    $r = $para.Range
    $xml = "<w:p $wns>$innerXml</w:p>"
    $r.InsertXML($xml)
}

# 1. "The system shall allow users to view their profile information."
#    -> "The " | [gramStart]"system shall"[gramEnd] | " allow users to view their profile information."
$p = Get-ParaByText "view their profile information"
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">The </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>system shall</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> allow users to view their profile information.</w:t></w:r>'
Set-ParaXml $p $inner

# 2. "The system shall allow users to update their profile information."
#    -> "The system " | [gramStart]"shall"[gramEnd] | " allow users to update their profile information."
$p = Get-ParaByText "update their profile information"
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">The system </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>shall</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> allow users to update their profile information.</w:t></w:r>'
Set-ParaXml $p $inner

# 3. "Constituents shall be able to view their representative's stats (...)"
#    -> "...view their " | [gramStart]"representative's"[gramEnd] | " stats (...)"
$p = Get-ParaByText "representative's stats"
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Constituents shall be able to view their </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>representative''s</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> stats (e.g., number of complaints, complaints resolved, pending complaints,</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> and</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> total meetups).</w:t></w:r>'
Set-ParaXml $p $inner

# 4. Collapse the 5 paragraphs from "...create and configure an assistant..."
#    through "The assistant shall be able to forward complaints to department."
#    (including the "Assistant:" Heading3 and its two bullets) into a single
#    "The Representative shall be able to forward complaints to the department."
#    bullet, re-using the first paragraph's pPr (ListParagraph / numId 24).
$pStart = Get-ParaByText "to create and configure"
$pEnd = Get-ParaByText "The assistant shall be able to forward complaints to department"
$full = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="24"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">The Representative shall be able to forward complaints to </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">the </w:t></w:r>' + `
    '<w:r><w:t>department.</w:t></w:r>'
$xml = "<w:p $wns>$inner</w:p>"
$full.InsertXML($xml)

# 5. "The system shall route uncategorized complaints to the representative and their assistant."
#    -> "The system shall route uncategorized complaints to the representative."
$d.Content.Find.Execute("The system shall route uncategorized complaints to the representative and their assistant.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The system shall route uncategorized complaints to the representative.", 2) | Out-Null

# 6 & 8. Move w:lastRenderedPageBreak from the "notify constituents" run to the
#        "Department:" heading run.
$p = Get-ParaByText "notify constituents about complaint status"
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="25"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>The system shall notify constituents about complaint status updates and upcoming virtual meetings.</w:t></w:r>'
Set-ParaXml $p $inner

# 7. "The system shall prioritize complaints containing sensitive keywords for urgent attention."
#    -> "The system " | [gramStart]"shall"[gramEnd] | " prioritize complaints containing sensitive keywords for urgent attention."
$p = Get-ParaByText "prioritize complaints containing sensitive"
$inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="25"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">The system </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>shall</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> prioritize complaints containing sensitive keywords for urgent attention.</w:t></w:r>'
Set-ParaXml $p $inner

# (InsertXML drops an explicit Heading3 w:pStyle on a from-scratch paragraph,
#  so rebuild without it, then reapply the paragraph style via .Style.)
$p = Get-ParaByText "^Department:"
$inner = '<w:r><w:lastRenderedPageBreak/><w:t>Department:</w:t></w:r>'
Set-ParaXml $p $inner
$p2 = Get-ParaByText "^Department:"
$p2.Style = "Heading 3"
